# Append two new rows (4 and 5) to sheet1, duplicating the existing
# data found in rows 2 and 3 (same venue/date/result/teams/batsman stats),
# exactly as described by the diff. All values in these rows are text
# (numberStoredAsText), so force a Text number format before assigning
# the values, which keeps Excel from re-interpreting numeric-looking
# strings ("0", "1", "2", "3", "0.00", "66.66") as real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = $ws.Range("A4:K5")
$newRows.NumberFormat = "@"

# Row 4 - duplicate of row 2 data
$ws.Cells.Item(4, 1).Value  = " Sharjah"
$ws.Cells.Item(4, 2).Value  = " October 26 2020"
$ws.Cells.Item(4, 3).Value  = "Kings XI won by 8 wickets (with 7 balls remaining)"
$ws.Cells.Item(4, 4).Value  = "Kolkata Knight Riders"
$ws.Cells.Item(4, 5).Value  = "Kings XI Punjab"
$ws.Cells.Item(4, 6).Value  = "Prasidh Krishna "
$ws.Cells.Item(4, 7).Value  = "0"
$ws.Cells.Item(4, 8).Value  = "1"
$ws.Cells.Item(4, 9).Value  = "0"
$ws.Cells.Item(4, 10).Value = "0"
$ws.Cells.Item(4, 11).Value = "0.00"

# Row 5 - duplicate of row 3 data
$ws.Cells.Item(5, 1).Value  = " Sharjah"
$ws.Cells.Item(5, 2).Value  = " October 12 2020"
$ws.Cells.Item(5, 3).Value  = "RCB won by 82 runs"
$ws.Cells.Item(5, 4).Value  = "Kolkata Knight Riders"
$ws.Cells.Item(5, 5).Value  = "Royal Challengers Bangalore"
$ws.Cells.Item(5, 6).Value  = "Prasidh Krishna "
$ws.Cells.Item(5, 7).Value  = "2"
$ws.Cells.Item(5, 8).Value  = "3"
$ws.Cells.Item(5, 9).Value  = "0"
$ws.Cells.Item(5, 10).Value = "0"
$ws.Cells.Item(5, 11).Value = "66.66"
